# NumerosSelecionados.xlsx — apply the "resultados.js" payout-id backfill +
# new selection rows described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small helper: write a value that must remain TEXT even when it is
# empty or all-digits (phone numbers / payment ids), which Excel would
# otherwise coerce to blank / a Number. A leading apostrophe (the normal
# "quote-prefix" convention) forces text in exactly those two cases;
# ordinary alphabetic text is left untouched so no extra formatting is
# introduced where it isn't needed.
function Set-TextValue {
    param($range, [string]$text)
    if ($text -eq "" -or $text -match '^[0-9]+$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# ---------------------------------------------------------------------
# 1) Backfill the previously-blank idPagamento (column D) for the four
#    rows that already existed.
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("D64") "78153849834"
Set-TextValue $ws.Range("D65") "77908618867"
Set-TextValue $ws.Range("D66") "78159597854"
Set-TextValue $ws.Range("D67") "77912975219"

# ---------------------------------------------------------------------
# 2) Append four new rows (70-73) produced by the new "10 numbers"
#    selection validation.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=70; A="Vitor Ito";                 B=1578424633; C="11988776655"; D=""; Nums=@(1,2,3,4,5,9,10,22,28,29);    O="Não" },
    @{ Row=71; A="Isabelly Silva Quintans";    B=7117522682; C="";           D=""; Nums=@(2,4,5,6,7,8,9,10,28,29);     O="Não" },
    @{ Row=72; A="Isabelly Silva Quintans";    B=7117522682; C="";           D=""; Nums=@(2,4,5,6,7,8,9,10,28,29);     O="Não" },
    @{ Row=73; A="Vitor Ito";                  B=1578424633; C="";           D=""; Nums=@(48,51,52,53,54,55,56,57,58,59); O="Não" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    Set-TextValue $ws.Cells.Item($row, 1) $r.A          # A - Nome
    $ws.Cells.Item($row, 2).Value = $r.B                # B - ID
    Set-TextValue $ws.Cells.Item($row, 3) $r.C          # C - Telefone
    Set-TextValue $ws.Cells.Item($row, 4) $r.D          # D - idPagamento

    for ($i = 0; $i -lt $r.Nums.Count; $i++) {
        $ws.Cells.Item($row, 5 + $i).Value = $r.Nums[$i]  # E..N - N1..N10
    }

    Set-TextValue $ws.Cells.Item($row, 15) $r.O         # O - Pagamento
}
